# Apply the commit's change: cell B4 on the (only) worksheet held the
# indicator title "16.6.1. Первичные расходы ..." (with a period after
# "16.6.1"). The edit drops that period -> "16.6.1 Первичные расходы ...".
#
# Editing the cell's text through the object model naturally causes Excel
# to rebuild the shared-strings table (dropping the old unique string and
# appending the edited one at the end), which is exactly the churn visible
# in the rest of the diff (all the si-index renumbering).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "16.6.1 Первичные расходы правительства в процентном отношении к первоначальному утвержденному бюджету в разбивке по секторам (по кодам бюджетной классификации или аналогичным категориям) "

# The author's edit also left the sheet's selection/active cell on B4
# (previously it was B2).
$ws.Range("B4").Select()

# Best-effort: the commit's workbook.xml bookViews entry also moved to a
# maximized-looking window (xWindow/yWindow = 0, larger windowWidth/
# windowHeight). Mirror that on the workbook window object.
$win = $wb.Windows.Item(1)
$win.WindowState = -4137
$win.Left = 0
$win.Top = 0
$win.Width = 28800
$win.Height = 11835
